$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the new "BASE data" sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "comparison with tool and data"

# Duplicate the first sheet (keeps headers/styles/column widths) and place it
# right after the original; this becomes the new "BASE data" sheet.
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "BASE data"

# --- Update the "BASE data" sheet's values (parsed via d3.nest as a nested
#     array rather than a map, so the re-computed/derived columns change) ---

# Row 3 (afg, governance delta 0)
$ws2.Range("F3").Value = 38.414000000000001
$ws2.Range("K3").Value = 45.871000000000002

# Row 4 (afg, governance delta 0.1)
$ws2.Range("F4").Value = 42.061
$ws2.Range("K4").Value = 47.341999999999999

# Row 5 (alb base row)
$ws2.Range("D5").Value = 1047.3499999999999

# Row 6 (alb, governance delta 0)
$ws2.Range("G6").Value = 64.546000000000006

# Row 7 (alb, governance delta 0.1)
$ws2.Range("G7").Value = 65.781000000000006

# Row 8 (alb base row, later years)
$ws2.Range("D8").Value = 1190.3599999999999

# Row 9 (alb, governance delta 0)
$ws2.Range("H9").Value = 97.846999999999994
$ws2.Range("I9").Value = 39.874000000000002

# Row 10 (alb, governance delta 0.1)
$ws2.Range("H10").Value = 98.36
$ws2.Range("I10").Value = 42.524000000000001

# Row 11 (ago base row)
$ws2.Range("D11").Value = 595.41
$ws2.Range("L11").Value = 91.94

# Row 12 (ago, governance delta 0)
$ws2.Range("L12").Value = 93.466999999999999
$ws2.Range("M12").Value = 99.881

# Row 13 (ago, governance delta 0.1)
$ws2.Range("L13").Value = 93.643000000000001
$ws2.Range("M13").Value = 99.893000000000001

# Remove the "tool value if different" column header + the now-unused J
# (schoolperc) values that don't apply to the BASE-data view; keep the
# N-column cells that still carry formatting but drop their leftover values.
$ws2.Range("N1").Clear() | Out-Null
$ws2.Range("J3").Clear() | Out-Null
$ws2.Range("J4").Clear() | Out-Null
$ws2.Range("N6").ClearContents() | Out-Null
$ws2.Range("N7").ClearContents() | Out-Null
$ws2.Range("N10").ClearContents() | Out-Null

# Column D (u5msurv) now holds wider values on this sheet, so it no longer
# keeps the narrow best-fit width it had on the comparison sheet.
$ws2.Columns.Item(4).ColumnWidth = 7.3

# --- Selections / active tab: BASE data is the sheet shown on reopen ---
$ws1.Activate()
$ws1.Range("C2").Select() | Out-Null

$ws2.Activate()
$ws2.Range("D2").Select() | Out-Null
